# Update "想去人数" (F column) figures across sheets, reflecting the
# newer snapshot of counts scraped for the gh-pages data output.

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1175
$ws1.Range("F3").Value  = 1991
$ws1.Range("F4").Value  = 631
$ws1.Range("F5").Value  = 1281
$ws1.Range("F7").Value  = 54
$ws1.Range("F8").Value  = 140
$ws1.Range("F9").Value  = 349
$ws1.Range("F10").Value = 135
$ws1.Range("F11").Value = 110
$ws1.Range("F12").Value = 877
$ws1.Range("F13").Value = 274
$ws1.Range("F14").Value = 140
$ws1.Range("F15").Value = 35
$ws1.Range("F16").Value = 118
$ws1.Range("F18").Value = 265
$ws1.Range("F19").Value = 712
$ws1.Range("F20").Value = 87
$ws1.Range("F22").Value = 212
$ws1.Range("F24").Value = 923
$ws1.Range("F26").Value = 204
$ws1.Range("F28").Value = 315
$ws1.Range("F30").Value = 24
$ws1.Range("F31").Value = 432

# Sheet: 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value  = 340
$ws2.Range("F11").Value = 134

# Sheet: 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 334

# Sheet: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 334
$ws4.Range("F3").Value  = 1175
$ws4.Range("F4").Value  = 1991
$ws4.Range("F5").Value  = 631
$ws4.Range("F6").Value  = 1281
$ws4.Range("F9").Value  = 54
$ws4.Range("F10").Value = 140
$ws4.Range("F11").Value = 349
$ws4.Range("F12").Value = 135
$ws4.Range("F13").Value = 110
$ws4.Range("F14").Value = 877
$ws4.Range("F15").Value = 274
$ws4.Range("F16").Value = 140
$ws4.Range("F18").Value = 35
$ws4.Range("F19").Value = 340
$ws4.Range("F20").Value = 118
$ws4.Range("F25").Value = 265
$ws4.Range("F26").Value = 712
$ws4.Range("F27").Value = 87
$ws4.Range("F29").Value = 212
$ws4.Range("F31").Value = 923
$ws4.Range("F35").Value = 204
$ws4.Range("F37").Value = 315
$ws4.Range("F39").Value = 134
$ws4.Range("F41").Value = 24
$ws4.Range("F43").Value = 432
